$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29
$ws.Range("A29").Value = 112526400
$ws.Range("B29").Value = 89551
$ws.Range("C29").Value = "Ovaliderad"
$ws.Range("D29").Value = "EN"
$ws.Range("E29").Value = 1110
$ws.Range("F29").Value = "Tallharticka"
$ws.Range("G29").Value = "Pelloporus triqueter"
$ws.Range("H29").Value = "(Pers.) Quél."
$ws.Range("P29").Value = "Stackmossen 500m Ö om, Sm"
$ws.Range("Q29").Value = 557833
$ws.Range("R29").Value = 6272869
$ws.Range("S29").Value = 10
$ws.Range("T29").Value = "Kalmar"
$ws.Range("U29").Value = "Kalmar"
$ws.Range("V29").Value = "Småland"
$ws.Range("W29").Value = "Mortorp"
$ws.Range("Y29").NumberFormat = "@"
$ws.Range("Y29").Value = "2023-10-04"
$ws.Range("Y29").Style = "Normal"
$ws.Range("AA29").NumberFormat = "@"
$ws.Range("AA29").Value = "2023-10-04"
$ws.Range("AA29").Style = "Normal"
$ws.Range("AC29").Value = "inga nya fruktkroppar"
$ws.Range("AD29").Value = $true
$ws.Range("AE29").Value = $false
$ws.Range("AG29").Value = $false
$ws.Range("AW29").Value = "Alexander Singer"
$ws.Range("AX29").Value = "Alexander Singer"

# Blank (but present) cells in row 29
$ws.Range("I29").Value = "​"
$ws.Range("J29").Value = "​"
$ws.Range("K29").Value = "​"
$ws.Range("N29").Value = "​"
$ws.Range("AF29").Value = "​"
$ws.Range("AT29").Value = "​"
$ws.Range("AY29").Value = "​"

# Row 30
$ws.Range("A30").Value = 112521407
$ws.Range("B30").Value = 90169
$ws.Range("C30").Value = "Ovaliderad"
$ws.Range("D30").Value = "LC"
$ws.Range("E30").Value = 6031
$ws.Range("F30").Value = "Blomkålssvamp"
$ws.Range("G30").Value = "Sparassis crispa"
$ws.Range("H30").Value = "(Wulfen:Fr.) Fr."
$ws.Range("P30").Value = "Stackmossen 500m Ö om, Sm"
$ws.Range("Q30").Value = 557833
$ws.Range("R30").Value = 6272869
$ws.Range("S30").Value = 10
$ws.Range("T30").Value = "Kalmar"
$ws.Range("U30").Value = "Kalmar"
$ws.Range("V30").Value = "Småland"
$ws.Range("W30").Value = "Mortorp"
$ws.Range("Y30").NumberFormat = "@"
$ws.Range("Y30").Value = "2023-10-04"
$ws.Range("Y30").Style = "Normal"
$ws.Range("AA30").NumberFormat = "@"
$ws.Range("AA30").Value = "2023-10-04"
$ws.Range("AA30").Style = "Normal"
$ws.Range("AD30").Value = $false
$ws.Range("AE30").Value = $false
$ws.Range("AG30").Value = $false
$ws.Range("AW30").Value = "Alexander Singer"
$ws.Range("AX30").Value = "Alexander Singer"

# Blank (but present) cells in row 30
$ws.Range("I30").Value = "​"
$ws.Range("K30").Value = "​"
$ws.Range("AT30").Value = "​"
$ws.Range("AY30").Value = "​"
